# The deck ships with two DrawingML theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" / "Office" colour scheme
#   ppt/theme/theme2.xml -> "Integral"     / "Red Violet" colour scheme
# (theme2.xml is the theme actually driving the slide master / slides;
#  theme1.xml is only linked from the notes master).
#
# The authored edit swaps the two themes' content wholesale: theme1.xml
# ends up holding the "Integral"/"Red Violet" definition and theme2.xml
# ends up holding the "Office Theme"/"Office" definition (font scheme and
# format scheme are already byte-identical between the two parts, so the
# only real content delta is the 12 colour-scheme slots plus the scheme
# names).
#
# Helper: pack R,G,B (0-255) into the little-endian integer that
# PowerPoint's COM automation uses for ColorFormat.RGB / the
# ThemeColorScheme item RGB property ( R + G*256 + B*65536 ).
function ColorVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# Target values: the "Office Theme" colour scheme (this is the scheme
# that theme1.xml currently has, and that theme2.xml must end up with).
# Order follows the standard MsoThemeColorSchemeIndex layout:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    (ColorVal 0x00 0x00 0x00),   # dk1      000000
    (ColorVal 0xFF 0xFF 0xFF),   # lt1      FFFFFF
    (ColorVal 0x44 0x54 0x6A),   # dk2      44546A
    (ColorVal 0xE7 0xE6 0xE6),   # lt2      E7E6E6
    (ColorVal 0x5B 0x9B 0xD5),   # accent1  5B9BD5
    (ColorVal 0xED 0x7D 0x31),   # accent2  ED7D31
    (ColorVal 0xA5 0xA5 0xA5),   # accent3  A5A5A5
    (ColorVal 0xFF 0xC0 0x00),   # accent4  FFC000
    (ColorVal 0x44 0x72 0xC4),   # accent5  4472C4
    (ColorVal 0x70 0xAD 0x47),   # accent6  70AD47
    (ColorVal 0x05 0x63 0xC1),   # hlink    0563C1
    (ColorVal 0x95 0x4F 0x72)    # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}

# Best-effort: also rename the design/theme/colour-scheme so the part
# reflects the "Office Theme" identity it is taking on.
try { $theme.Name = "Office Theme" } catch {}
try { $colorScheme.Name = "Office" } catch {}
try {
    $design = $p.Designs.Item(1)
    $design.Name = "Office Theme"
} catch {}
